$wb = $excel.ActiveWorkbook

# Sheet 1
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(98, 8).Value = 2179.4167
$ws.Cells.Item(98, 9).Value = 1239.3334
$ws.Cells.Item(98, 11).Value = 1239.3334
$ws.Cells.Item(98, 13).Value = 258.6666
$ws.Cells.Item(122, 8).Value = 2179.4167
$ws.Cells.Item(122, 9).Value = 1239.3334
$ws.Cells.Item(122, 11).Value = 3718.0002
$ws.Cells.Item(122, 13).Value = -1268.0002
$ws.Cells.Item(125, 8).Value = 7379.8
$ws.Cells.Item(125, 9).Value = 3999
$ws.Cells.Item(125, 11).Value = 35991
$ws.Cells.Item(125, 13).Value = -33531
$ws.Cells.Item(132, 8).Value = 1600.8605
$ws.Cells.Item(132, 9).Value = 1444.8049
$ws.Cells.Item(132, 10).Value = 4800
$ws.Cells.Item(132, 11).Value = 4334.4147
$ws.Cells.Item(132, 12).Value = 14400
$ws.Cells.Item(132, 13).Value = -1804.4147
$ws.Cells.Item(132, 14).Value = -19460
$ws.Cells.Item(137, 8).Value = 1985503
$ws.Cells.Item(137, 9).Value = 1127.625
$ws.Cells.Item(137, 10).Value = 4631337
$ws.Cells.Item(137, 11).Value = 3382.875
$ws.Cells.Item(137, 12).Value = 13894011
$ws.Cells.Item(137, 13).Value = -832.875
$ws.Cells.Item(137, 14).Value = -13899111
$ws.Cells.Item(138, 8).Value = 2470.75
$ws.Cells.Item(138, 10).Value = 2750
$ws.Cells.Item(138, 12).Value = 8250
$ws.Cells.Item(138, 14).Value = -18530

# Sheet 2
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(2, 8).Value = 2370.75
$ws.Cells.Item(2, 10).Value = 3012.25
$ws.Cells.Item(2, 12).Value = 3012.25
$ws.Cells.Item(2, 14).Value = -3238.25
$ws.Cells.Item(4, 8).Value = 190.16667
$ws.Cells.Item(4, 9).Value = 200.2
$ws.Cells.Item(4, 11).Value = 200.2
$ws.Cells.Item(4, 13).Value = -84.19999999999999
$ws.Cells.Item(45, 8).Value = 5699.6665
$ws.Cells.Item(45, 9).Value = 5599.5
$ws.Cells.Item(45, 11).Value = 5599.5
$ws.Cells.Item(45, 13).Value = -5222.5
$ws.Cells.Item(63, 8).Value = 4339.6
$ws.Cells.Item(66, 8).Value = 4339.6
$ws.Cells.Item(74, 8).Value = 2758.5
$ws.Cells.Item(74, 9).Value = 2693.25
$ws.Cells.Item(74, 10).Value = 3150
$ws.Cells.Item(74, 11).Value = 2693.25
$ws.Cells.Item(74, 12).Value = 3150
$ws.Cells.Item(74, 13).Value = -1819.25
$ws.Cells.Item(74, 14).Value = -4898
$ws.Cells.Item(77, 8).Value = 2758.5
$ws.Cells.Item(77, 9).Value = 2693.25
$ws.Cells.Item(77, 10).Value = 3150
$ws.Cells.Item(77, 11).Value = 13466.25
$ws.Cells.Item(77, 12).Value = 15750
$ws.Cells.Item(77, 13).Value = -9098.25
$ws.Cells.Item(77, 14).Value = -24486
$ws.Cells.Item(97, 8).Value = 3422.2307
$ws.Cells.Item(97, 9).Value = 2029.8334
$ws.Cells.Item(97, 11).Value = 2029.8334
$ws.Cells.Item(97, 13).Value = -1533.8334
$ws.Cells.Item(116, 8).Value = 2370.75
$ws.Cells.Item(116, 10).Value = 3012.25
$ws.Cells.Item(116, 12).Value = 3012.25
$ws.Cells.Item(116, 14).Value = -7600.25

# Sheet 3
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(3, 8).Value = 2370.75
$ws.Cells.Item(3, 10).Value = 3012.25
$ws.Cells.Item(3, 12).Value = 3012.25
$ws.Cells.Item(3, 14).Value = -3240.25
$ws.Cells.Item(26, 8).Value = 57302.6
$ws.Cells.Item(26, 9).Value = 25504.5
$ws.Cells.Item(26, 11).Value = 25504.5
$ws.Cells.Item(26, 13).Value = -25212.5
$ws.Cells.Item(42, 8).Value = 288011.5
$ws.Cells.Item(42, 10).Value = 288011.5
$ws.Cells.Item(42, 12).Value = 288011.5
$ws.Cells.Item(42, 14).Value = -288667.5
$ws.Cells.Item(96, 8).Value = 26942.166
$ws.Cells.Item(96, 9).Value = 20413.75
$ws.Cells.Item(96, 10).Value = 39999
$ws.Cells.Item(96, 11).Value = 20413.75
$ws.Cells.Item(96, 12).Value = 39999
$ws.Cells.Item(96, 13).Value = -17667.75
$ws.Cells.Item(96, 14).Value = -45491
$ws.Cells.Item(125, 8).Value = 121000
$ws.Cells.Item(125, 10).Value = 121000
$ws.Cells.Item(125, 12).Value = 121000
$ws.Cells.Item(125, 14).Value = -130840

# Sheet 4
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(31, 8).Value = 4236.6387
$ws.Cells.Item(31, 9).Value = 1849.5
$ws.Cells.Item(31, 10).Value = 4714.067
$ws.Cells.Item(31, 11).Value = 1849.5
$ws.Cells.Item(31, 12).Value = 4714.067
$ws.Cells.Item(31, 13).Value = -1554.5
$ws.Cells.Item(31, 14).Value = -5304.067
$ws.Cells.Item(34, 8).Value = 4236.6387
$ws.Cells.Item(34, 9).Value = 1849.5
$ws.Cells.Item(34, 10).Value = 4714.067
$ws.Cells.Item(34, 11).Value = 1849.5
$ws.Cells.Item(34, 12).Value = 4714.067
$ws.Cells.Item(34, 13).Value = -1647.5
$ws.Cells.Item(34, 14).Value = -5118.067
$ws.Cells.Item(99, 8).Value = 2256
$ws.Cells.Item(99, 9).Value = 2256
$ws.Cells.Item(99, 11).Value = 2256
$ws.Cells.Item(99, 13).Value = -758
$ws.Cells.Item(122, 8).Value = 3607.8293
$ws.Cells.Item(122, 10).Value = 4726.143
$ws.Cells.Item(122, 12).Value = 14178.429
$ws.Cells.Item(122, 14).Value = -19078.429
$ws.Cells.Item(126, 8).Value = 2256
$ws.Cells.Item(126, 9).Value = 2256
$ws.Cells.Item(126, 11).Value = 6768
$ws.Cells.Item(126, 13).Value = -4298
$ws.Cells.Item(132, 8).Value = 2842.925
$ws.Cells.Item(132, 9).Value = 2459.0303
$ws.Cells.Item(132, 10).Value = 4652.7144
$ws.Cells.Item(132, 11).Value = 7377.090899999999
$ws.Cells.Item(132, 12).Value = 13958.1432
$ws.Cells.Item(132, 13).Value = -4847.090899999999
$ws.Cells.Item(132, 14).Value = -19018.1432
$ws.Cells.Item(141, 8).Value = 680998.7
$ws.Cells.Item(141, 10).Value = 1952998
$ws.Cells.Item(141, 12).Value = 1952998
$ws.Cells.Item(141, 14).Value = -1963358

# Sheet 5
$ws = $wb.Worksheets.Item(5)
$ws.Cells.Item(60, 8).Value = 1003.6667
$ws.Cells.Item(60, 10).Value = 1003.6667
$ws.Cells.Item(60, 12).Value = 3011.0001
$ws.Cells.Item(60, 14).Value = -3513.0001
$ws.Cells.Item(117, 8).Value = 3368.3333
$ws.Cells.Item(117, 9).Value = 300
$ws.Cells.Item(117, 10).Value = 3751.875
$ws.Cells.Item(117, 11).Value = 900
$ws.Cells.Item(117, 12).Value = 11255.625
$ws.Cells.Item(117, 13).Value = 2542
$ws.Cells.Item(117, 14).Value = -18139.625
$ws.Cells.Item(131, 8).Value = 1750.5358
$ws.Cells.Item(131, 10).Value = 1815.4783
$ws.Cells.Item(131, 12).Value = 5446.4349
$ws.Cells.Item(131, 14).Value = -15526.4349
$ws.Cells.Item(132, 8).Value = 478659.72
$ws.Cells.Item(132, 9).Value = 1826.6
$ws.Cells.Item(132, 11).Value = 16439.4
$ws.Cells.Item(132, 13).Value = -13909.4
$ws.Cells.Item(137, 8).Value = 2369.8462
$ws.Cells.Item(137, 9).Value = 1880.8
$ws.Cells.Item(137, 10).Value = 4000
$ws.Cells.Item(137, 11).Value = 5642.4
$ws.Cells.Item(137, 12).Value = 12000
$ws.Cells.Item(137, 13).Value = -542.3999999999996
$ws.Cells.Item(137, 14).Value = -22200

# Sheet 6
$ws = $wb.Worksheets.Item(6)
$ws.Cells.Item(34, 8).Value = 31250
$ws.Cells.Item(76, 8).Value = 31250
$ws.Cells.Item(79, 8).Value = 31250
$ws.Cells.Item(113, 8).Value = 14621.0625
$ws.Cells.Item(113, 9).Value = 17276.572
$ws.Cells.Item(113, 10).Value = 12555.667
$ws.Cells.Item(113, 11).Value = 17276.572
$ws.Cells.Item(113, 12).Value = 12555.667
$ws.Cells.Item(113, 13).Value = -15106.572
$ws.Cells.Item(113, 14).Value = -16895.667

# Sheet 7
$ws = $wb.Worksheets.Item(7)
$ws.Cells.Item(25, 8).Value = 0
$ws.Cells.Item(25, 9).Value = 0
$ws.Cells.Item(25, 11).Value = 0
$ws.Cells.Item(40, 8).Value = 2975.524
$ws.Cells.Item(40, 9).Value = 2102.4
$ws.Cells.Item(40, 11).Value = 2102.4
$ws.Cells.Item(40, 13).Value = -1966.4
$ws.Cells.Item(46, 8).Value = 10594.625
$ws.Cells.Item(46, 9).Value = 0
$ws.Cells.Item(46, 10).Value = 10594.625
$ws.Cells.Item(46, 11).Value = 0
$ws.Cells.Item(46, 12).Value = 10594.625
$ws.Cells.Item(46, 14).Value = -10970.625
$ws.Cells.Item(61, 8).Value = 1540.2632
$ws.Cells.Item(61, 9).Value = 1456.7333
$ws.Cells.Item(61, 10).Value = 1853.5
$ws.Cells.Item(61, 11).Value = 1456.7333
$ws.Cells.Item(61, 12).Value = 1853.5
$ws.Cells.Item(61, 13).Value = -1254.7333
$ws.Cells.Item(61, 14).Value = -2257.5
$ws.Cells.Item(80, 8).Value = 130998.664
$ws.Cells.Item(80, 9).Value = 116998
$ws.Cells.Item(80, 10).Value = 159000
$ws.Cells.Item(80, 11).Value = 116998
$ws.Cells.Item(80, 12).Value = 159000
$ws.Cells.Item(80, 13).Value = -115875
$ws.Cells.Item(80, 14).Value = -161246
$ws.Cells.Item(83, 8).Value = 130998.664
$ws.Cells.Item(83, 9).Value = 116998
$ws.Cells.Item(83, 10).Value = 159000
$ws.Cells.Item(83, 11).Value = 350994
$ws.Cells.Item(83, 12).Value = 477000
$ws.Cells.Item(83, 13).Value = -345378
$ws.Cells.Item(83, 14).Value = -488232
$ws.Cells.Item(87, 8).Value = 96999.5
$ws.Cells.Item(87, 9).Value = 80000
$ws.Cells.Item(87, 10).Value = 113999
$ws.Cells.Item(87, 11).Value = 80000
$ws.Cells.Item(87, 12).Value = 113999
$ws.Cells.Item(87, 13).Value = -78877
$ws.Cells.Item(87, 14).Value = -116245
$ws.Cells.Item(90, 8).Value = 96999.5
$ws.Cells.Item(90, 9).Value = 80000
$ws.Cells.Item(90, 10).Value = 113999
$ws.Cells.Item(90, 11).Value = 240000
$ws.Cells.Item(90, 12).Value = 341997
$ws.Cells.Item(90, 13).Value = -234384
$ws.Cells.Item(90, 14).Value = -353229
$ws.Cells.Item(113, 8).Value = 1540.2632
$ws.Cells.Item(113, 9).Value = 1456.7333
$ws.Cells.Item(113, 10).Value = 1853.5
$ws.Cells.Item(113, 11).Value = 1456.7333
$ws.Cells.Item(113, 12).Value = 1853.5
$ws.Cells.Item(113, 13).Value = 713.2666999999999
$ws.Cells.Item(113, 14).Value = -6193.5
$ws.Cells.Item(132, 8).Value = 5466.1904
$ws.Cells.Item(132, 9).Value = 5206.375
$ws.Cells.Item(132, 11).Value = 15619.125
$ws.Cells.Item(132, 13).Value = -13089.125
$ws.Cells.Item(25, 13).ClearContents()
$ws.Cells.Item(46, 13).ClearContents()

# Sheet 8
$ws = $wb.Worksheets.Item(8)
$ws.Cells.Item(96, 8).Value = 8434.647000000001
$ws.Cells.Item(96, 9).Value = 5047
$ws.Cells.Item(96, 10).Value = 11445.889
$ws.Cells.Item(96, 11).Value = 5047
$ws.Cells.Item(96, 12).Value = 11445.889
$ws.Cells.Item(96, 13).Value = -3674
$ws.Cells.Item(96, 14).Value = -14191.889
$ws.Cells.Item(126, 8).Value = 7557.4287
$ws.Cells.Item(126, 9).Value = 9700.75
$ws.Cells.Item(126, 11).Value = 29102.25
$ws.Cells.Item(126, 13).Value = -26632.25
$ws.Cells.Item(136, 8).Value = 15464374
$ws.Cells.Item(136, 9).Value = 2602.3333
$ws.Cells.Item(136, 11).Value = 7806.999899999999
$ws.Cells.Item(136, 13).Value = -5256.999899999999
